# Applies the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '53.513.33'
$ws.Range('E2').Value = '  -4.44%  '
$ws.Range('D3').Value = '2.191.26'
$ws.Range('E3').Value = '  -7.21%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '483.34'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '124.34'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.85%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('E8').Value = '  -4.60%  '
$ws.Range('D9').Value = '2.208.26'
$ws.Range('E9').Value = '  -6.56%  '
$ws.Range('E10').Value = '  -6.75%  '
$ws.Range('E11').Value = '  -1.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.61'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.90%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.313'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.17%  '
$ws.Range('D14').Value = '2.581.58'
$ws.Range('E14').Value = '  -7.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.07'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.24%  '
$ws.Range('D16').Value = '53.422.80'
$ws.Range('E16').Value = '  -4.61%  '
$ws.Range('E17').Value = '  -3.64%  '
$ws.Range('D18').Value = '2.201.65'
$ws.Range('E18').Value = '  -4.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.54'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.61%  '
$ws.Range('E20').Value = '  -2.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '294.08'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.08%  '
$ws.Range('E22').Value = '  -3.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.49'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.996'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.365'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.08%  '
$ws.Range('D27').Value = '2.295.42'
$ws.Range('E27').Value = '  -7.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.145'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.20%  '
$ws.Range('E29').Value = '  -3.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '165.75'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.33%  '
$ws.Range('E31').Value = '  -4.06%  '
$ws.Range('E32').Value = '  -0.19%  '
$ws.Range('E33').Value = '  -0.34%  '
$ws.Range('D34').Value = '0.0₃0656'
$ws.Range('E34').Value = '  -7.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.65'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.63%  '
$ws.Range('E36').Value = '  -2.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.26'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.73%  '
$ws.Range('E38').Value = '  -2.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.821'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.86%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '35.70'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.08%  '
$ws.Range('B41').Value = 'NEARProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.52'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.367'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.77%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.35'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.44%  '
$ws.Range('E44').Value = '  -3.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '123.23'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.72'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.36%  '
$ws.Range('E47').Value = '  -3.11%  '
$ws.Range('E48').Value = '  -5.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '228.38'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0467'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.86%  '
$ws.Range('E51').Value = '  -3.80%  '
